# Update cryptos list prices and volume(1h) percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.563.30"
$ws.Range("D3").Value = "3.148.94"
$ws.Range("E3").Value = "  -4.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'526.42"
$ws.Range("E5").Value = "  -4.89%  "
$ws.Range("D6").Value = "'135.39"
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.148.63"
$ws.Range("E8").Value = "  -4.09%  "
$ws.Range("D9").Value = "'0.443"
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("D10").Value = "'7.23"
$ws.Range("E10").Value = "  -6.76%  "
$ws.Range("E11").Value = "  -8.07%  "
$ws.Range("E12").Value = "  -6.20%  "
$ws.Range("D13").Value = "3.687.49"
$ws.Range("E13").Value = "  -4.14%  "
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "'25.62"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D16").Value = "3.146.32"
$ws.Range("E16").Value = "  -4.03%  "
$ws.Range("D17").Value = "58.500.54"
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("E18").Value = "  -6.17%  "
$ws.Range("D19").Value = "'5.78"
$ws.Range("D20").Value = "'13.05"
$ws.Range("E20").Value = "  -5.41%  "
$ws.Range("D21").Value = "'7.95"
$ws.Range("E21").Value = "  -6.87%  "
$ws.Range("D22").Value = "'345.14"
$ws.Range("E22").Value = "  -7.23%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -4.41%  "
$ws.Range("D25").Value = "'68.30"
$ws.Range("E25").Value = "  -7.07%  "
$ws.Range("D26").Value = "3.273.24"
$ws.Range("E26").Value = "  -4.29%  "
$ws.Range("D27").Value = "'0.171"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "0.0₃0959"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "'6.83"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -7.47%  "
$ws.Range("D33").Value = "'6.91"
$ws.Range("E33").Value = "  -7.46%  "
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("E36").Value = "  -4.43%  "
$ws.Range("D37").Value = "'157.32"
$ws.Range("E37").Value = "  -5.37%  "
$ws.Range("E38").Value = "  -5.74%  "
$ws.Range("D40").Value = "'0.0687"
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("D41").Value = "3.179.78"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("D42").Value = "'24.30"
$ws.Range("E42").Value = "  -6.98%  "
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "'0.694"
$ws.Range("E45").Value = "  -6.94%  "
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  -7.63%  "
$ws.Range("D49").Value = "2.275.71"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "'6.21"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'20.86"
$ws.Range("E51").Value = "  -1.58%  "

